# Add NL Wave 8 data (row 17) to the "NL" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NL")

# New data row describing the NL wave 8 survey file.
$ws.Cells.Item(17, 1).Value = 6                     # A17 survey_version
$ws.Cells.Item(17, 2).Value = 0                      # B17 locked
$ws.Cells.Item(17, 3).Value = "nl"                   # C17 country
$ws.Cells.Item(17, 4).Value = 16                     # D17 survey_round
$ws.Cells.Item(17, 5).Value = "B"                    # E17 panel
$ws.Cells.Item(17, 6).Value = 8                      # F17 wave

# G17 date_recieved - copy the date formatting used by the row above,
# then set the new date's underlying serial value (2021-04-09).
$ws.Cells.Item(16, 7).Copy($ws.Cells.Item(17, 7))
$ws.Cells.Item(17, 7).Value = 44295

$ws.Cells.Item(17, 8).Value = "20-090916_NL_Wave8_Final_v1_090421_IntClientUse"  # H17 spss_name

# I17 r_name - same formula pattern used by the other rows in the column.
$ws.Cells.Item(17, 9).Formula = "=C17&""_""&""sr""&TEXT(D17,""00"")&""_""&YEAR(G17)&TEXT(G17,""MM"")&TEXT(G17,""DD"")&""_p""&E17&""_wv""&TEXT(F17,""00"")&"""""

# Update the active selection to match the edited workbook.
$ws.Activate()
$ws.Range("B17").Select()
